$d = $word.ActiveDocument
$d.Content.Find.Execute("[if player is not Bosmer] ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[if player is not Bosmer or Altmer] ", 2)
